# KSA_Cities.xlsx update:
#  - Rename the existing "Tabuk Principal" entry (row 200) to "Tabuk Province Emirate"
#  - Append three new city rows (Al Baha Province Emirate, Al Houta, Raqai)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 200: "Tabuk Principal" -> "Tabuk Province Emirate" -------------
$ws.Range("A200").Value = "Tabuk Province Emirate"
$ws.Range("B200").Value = "Tabuk Province Emirate"
# C200 (Arabic name), D200/E200 (lat/long) and F200/G200 (area/region) stay as-is.

# --- Duplicate the formatting of row 200 down into the three new rows ---------
$ws.Range("A200:G200").Copy()
$ws.Range("A201:G203").PasteSpecial(-4122)

# --- Row 201: Al Baha Province Emirate -----------------------------------------
$ws.Range("A201").Value = "Al Baha Province Emirate"
$ws.Range("B201").Value = "Al Baha Province Emirate"
$ws.Range("C201").Value = "الباحة (مقرالامارة)"
$ws.Range("D201").Value = 20.01211
$ws.Range("E201").Value = 41.467326999999997
$ws.Range("F201").Value = "منطقة الباحة"
$ws.Range("G201").Value = "جنوب المملكة"

# --- Row 202: Al Houta -----------------------------------------------------------
$ws.Range("A202").Value = "Al Houta"
$ws.Range("B202").Value = "Al Houta"
$ws.Range("C202").Value = "الحوطة"
$ws.Range("D202").Value = 23.525065000000001
$ws.Range("E202").Value = 46.845830999999997
$ws.Range("F202").Value = "منطقة الرياض"
$ws.Range("G202").Value = "وسط المملكة"

# --- Row 203: Raqai ---------------------------------------------------------------
$ws.Range("A203").Value = "Raqai"
$ws.Range("B203").Value = "Raqai"
$ws.Range("C203").Value = "الرقعي"
$ws.Range("D203").Value = 29.066296999999999
$ws.Range("E203").Value = 46.638584000000002
$ws.Range("F203").Value = "المنطقة الشرقية"
$ws.Range("G203").Value = "شرق المملكة"

# --- Refresh the full-table selection to cover the new rows --------------------
$ws.Range("A1:G203").Select()
